# Error corrections in the "Data" sheet's Value column (G) - re-calculated
# reporting data. Each target value already exists elsewhere in the shared
# string table, so we only need to repoint the affected cells to the
# corrected text. The value is text (matching the existing "t=s" cells in
# column G), so it is written with a leading apostrophe to stop Excel from
# re-interpreting the numeric-looking text as a number; the style is then
# reset to "Normal" so no stray per-cell formatting (e.g. quote-prefix) is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

function Set-TextValue($address, $text) {
    $ws.Range($address).Value = "'" + $text
    $ws.Range($address).Style = "Normal"
}

Set-TextValue "G8"   "2"
Set-TextValue "G11"  "1"
Set-TextValue "G14"  "2"
Set-TextValue "G15"  "5"
Set-TextValue "G17"  "1"
Set-TextValue "G30"  "5"
Set-TextValue "G31"  "10"
Set-TextValue "G33"  "30"
Set-TextValue "G34"  "1"
Set-TextValue "G291" "4"
Set-TextValue "G294" "1"
Set-TextValue "G297" "5"
Set-TextValue "G298" "2"
Set-TextValue "G300" "2"
Set-TextValue "G313" "18"
Set-TextValue "G314" "13"
Set-TextValue "G317" "2"
Set-TextValue "G574" "4"
Set-TextValue "G577" "1"
Set-TextValue "G580" "5"
Set-TextValue "G581" "2"
Set-TextValue "G583" "2"
Set-TextValue "G596" "18"
Set-TextValue "G597" "13"
Set-TextValue "G600" "2"
Set-TextValue "G857" "5"
Set-TextValue "G860" "1"
Set-TextValue "G863" "6"
Set-TextValue "G864" "5"
Set-TextValue "G866" "2"
Set-TextValue "G879" "27"
Set-TextValue "G880" "26"
Set-TextValue "G883" "7"
